$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 131, shifting existing rows 131-157 down to 132-158.
$ws.Rows("131").Insert()

# Populate the newly inserted row 131 with the new weekly price record.
$ws.Range("A131").Value = 8
$ws.Range("B131").Value = "Terminal La Palmera de La Serena"
$ws.Range("C131").Value = "Coquimbo"
$ws.Range("D131").Value = 45211
$ws.Range("E131").Value = 4
$ws.Range("F131").Value = 100112028
$ws.Range("G131").Value = "Sandia"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 2000
$ws.Range("K131").Value = 800
$ws.Range("L131").Value = 900
$ws.Range("M131").Value = 850
$ws.Range("N131").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O131").Value = "Perú"
$ws.Range("P131").Value = 850
$ws.Range("Q131").Value = 1
$ws.Range("R131").Value = "Hortaliza"
